# Applies the 'Week 7' data refresh to the monthly chart workbook.
# Re-populates each report sheet (track/artist/label x YouTube/1001Tracklists/Soundcloud)
# with the updated weekly figures, preserving header formatting and the bordered/bold
# style used on the key column of the artist/label summary sheets.
$wb = $excel.ActiveWorkbook

function Set-SheetData {
    param($SheetName, $Rows, $NumCols, $StyleColA, $OrigDataRows)
    $ws = $wb.Worksheets.Item($SheetName)
    $styleSource = $null
    if ($StyleColA) {
        $styleSource = $ws.Range("A2")
    }
    for ($i = 0; $i -lt $Rows.Count; $i++) {
        $r = 2 + $i
        $rowData = $Rows[$i]
        if ($StyleColA -and ($i -ge $OrigDataRows)) {
            $styleSource.Copy()
            $ws.Cells.Item($r, 1).PasteSpecial(-4122)
        }
        for ($c = 0; $c -lt $NumCols; $c++) {
            $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
        }
    }
}

# By_Track_YouTube
$rows = @(
    ,@("Blasterjaxx, Dr Phunk", "Here Without You", "Maxximize Records, Spinnin' Records", 164407)
    ,@("Bad Computer", "Destroy Me", "Monstercat", 117028)
    ,@("Koven", "Light Up", "Monstercat", 87606)
    ,@("Ellis", "Orbit", "Monstercat", 77004)
    ,@("Centineo, Niles Mason", "Fearless", "NONE", 66444)
    ,@("Masayoshi Iimori", "In My Soul", "Monstercat", 52535)
    ,@("SWACQ", "Purification", "HEXAGON", 25247)
    ,@("KAAZE", "Midnight Runners", "Revealed Music", 21714)
    ,@("Trilane, Charlie Ray, Jordan Grace", "In Heaven", "Protocol Recordings", 7277)
    ,@("Nokturn", "Stargazing", "NONE", 311)
)
Set-SheetData "By_Track_YouTube" $rows 4 $false 7

# By_Track_1001Tracklists
$rows = @(
    ,@("SWACQ", "Purification", "HEXAGON", 42, 48)
    ,@("KAAZE", "Midnight Runners", "Revealed Music", 25, 27)
    ,@("Trilane, Charlie Ray, Jordan Grace", "In Heaven", "Protocol Recordings", 12, 13)
    ,@("Blasterjaxx, Dr Phunk", "Here Without You", "Maxximize Records, Spinnin' Records", 10, 30)
    ,@("Ellis", "Orbit", "Monstercat", 5, 5)
    ,@("Bad Computer", "Destroy Me", "Monstercat", 3, 4)
    ,@("Koven", "Light Up", "Monstercat", 2, 3)
    ,@("Masayoshi Iimori", "In My Soul", "Monstercat", 1, 3)
    ,@("Centineo, Niles Mason", "Fearless", "NONE", 0, 0)
    ,@("Nokturn", "Stargazing", "NONE", 0, 0)
)
Set-SheetData "By_Track_1001Tracklists" $rows 5 $false 7

# By_Track_Soundcloud
$rows = @(
    ,@("Bad Computer", "Destroy Me", "Monstercat", 39642)
    ,@("Ellis", "Orbit", "Monstercat", 21256)
    ,@("Masayoshi Iimori", "In My Soul", "Monstercat", 13813)
    ,@("Koven", "Light Up", "Monstercat", 13795)
    ,@("Blasterjaxx, Dr Phunk", "Here Without You", "Maxximize Records, Spinnin' Records", 11363)
    ,@("SWACQ", "Purification", "HEXAGON", 7750)
    ,@("KAAZE", "Midnight Runners", "Revealed Music", 4350)
    ,@("Centineo, Niles Mason", "Fearless", "NONE", 2993)
    ,@("Trilane, Charlie Ray, Jordan Grace", "In Heaven", "Protocol Recordings", 1342)
    ,@("Nokturn", "Stargazing", "NONE", 35)
)
Set-SheetData "By_Track_Soundcloud" $rows 4 $false 7

# By_Artist_YouTube
$rows = @(
    ,@("Blasterjaxx", 164407)
    ,@("Dr Phunk", 164407)
    ,@("Bad Computer", 117028)
    ,@("Koven", 87606)
    ,@("Ellis", 77004)
    ,@("Centineo", 66444)
    ,@("Niles Mason", 66444)
    ,@("Masayoshi Iimori", 52535)
    ,@("SWACQ", 25247)
    ,@("KAAZE", 21714)
    ,@("Charlie Ray", 7277)
    ,@("Jordan Grace", 7277)
    ,@("Trilane", 7277)
    ,@("Nokturn", 311)
)
Set-SheetData "By_Artist_YouTube" $rows 2 $true 8

# By_Artist_1001Tracklists
$rows = @(
    ,@("SWACQ", 42, 48)
    ,@("KAAZE", 25, 27)
    ,@("Charlie Ray", 12, 13)
    ,@("Jordan Grace", 12, 13)
    ,@("Trilane", 12, 13)
    ,@("Blasterjaxx", 10, 30)
    ,@("Dr Phunk", 10, 30)
    ,@("Ellis", 5, 5)
    ,@("Bad Computer", 3, 4)
    ,@("Koven", 2, 3)
    ,@("Masayoshi Iimori", 1, 3)
    ,@("Centineo", 0, 0)
    ,@("Niles Mason", 0, 0)
    ,@("Nokturn", 0, 0)
)
Set-SheetData "By_Artist_1001Tracklists" $rows 3 $true 8

# By_Artist_Soundcloud
$rows = @(
    ,@("Bad Computer", 39642)
    ,@("Ellis", 21256)
    ,@("Masayoshi Iimori", 13813)
    ,@("Koven", 13795)
    ,@("Blasterjaxx", 11363)
    ,@("Dr Phunk", 11363)
    ,@("SWACQ", 7750)
    ,@("KAAZE", 4350)
    ,@("Centineo", 2993)
    ,@("Niles Mason", 2993)
    ,@("Charlie Ray", 1342)
    ,@("Jordan Grace", 1342)
    ,@("Trilane", 1342)
    ,@("Nokturn", 35)
)
Set-SheetData "By_Artist_Soundcloud" $rows 2 $true 8

# By_Label_YouTube
$rows = @(
    ,@("Monstercat", 334173)
    ,@("Maxximize Records", 164407)
    ,@("Spinnin' Records", 164407)
    ,@("HEXAGON", 25247)
    ,@("Revealed Music", 21714)
    ,@("Protocol Recordings", 7277)
)
Set-SheetData "By_Label_YouTube" $rows 2 $true 3

# By_Label_1001Tracklists
$rows = @(
    ,@("HEXAGON", 42, 48)
    ,@("Revealed Music", 25, 27)
    ,@("Protocol Recordings", 12, 13)
    ,@("Monstercat", 11, 15)
    ,@("Maxximize Records", 10, 30)
    ,@("Spinnin' Records", 10, 30)
)
Set-SheetData "By_Label_1001Tracklists" $rows 3 $true 3

# By_Label_Soundcloud
$rows = @(
    ,@("Monstercat", 88506)
    ,@("Maxximize Records", 11363)
    ,@("Spinnin' Records", 11363)
    ,@("HEXAGON", 7750)
    ,@("Revealed Music", 4350)
    ,@("Protocol Recordings", 1342)
)
Set-SheetData "By_Label_Soundcloud" $rows 2 $true 3
